# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#   populated with the per-fund holdings for the new quarter.
# - Update the "总计" (totals) sheet with a new leading row summarizing the
#   2022-Q1 quarter, shifting the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force a (possibly numeric-looking) string to be stored as literal
    # text rather than letting Excel auto-convert it to a number, while
    # not leaving behind any custom number-format/style on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4" (i.e. right
#    before "总计", which keeps it in chronological order).
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Clone the header/formatting of an existing per-quarter sheet (same
# column layout & cell styles) rather than building styles from scratch.
$q4.Range("B1:H2").Copy($newSheet.Range("B1:H2"))
$q4.Range("A2").Copy($newSheet.Range("A2"))
$q4.Range("A2").Copy($newSheet.Range("A3"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 009327 东兴兴晟混合A
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "009327"
$newSheet.Range("C2").Value = "东兴兴晟混合A"
Set-TextValue $newSheet.Range("D2") "0.34"
Set-TextValue $newSheet.Range("E2") "79.83"
Set-TextValue $newSheet.Range("F2") "1.04"
Set-TextValue $newSheet.Range("G2") "0.0035"
$newSheet.Range("H2").Value = 4

# Row 3 - 009328 东兴兴晟混合C
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "009328"
$newSheet.Range("C3").Value = "东兴兴晟混合C"
Set-TextValue $newSheet.Range("D3") "0.08"
Set-TextValue $newSheet.Range("E3") "79.83"
Set-TextValue $newSheet.Range("F3") "1.04"
Set-TextValue $newSheet.Range("G3") "0.0008"
$newSheet.Range("H3").Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: shift existing quarter rows down one row
#    and insert the 2022-Q1 summary at the top of the data (row 2).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$row2B = $total.Range("B2").Value2
$row2C = $total.Range("C2").Value2
$row2D = $total.Range("D2").Value2
$row3B = $total.Range("B3").Value2
$row3C = $total.Range("C3").Value2
$row3D = $total.Range("D3").Value2
$row4B = $total.Range("B4").Value2
$row4C = $total.Range("C4").Value2
$row4D = $total.Range("D4").Value2

# Row 5 is brand new - clone the index-column style (border/bold) from
# an existing row in column A before writing its value.
$total.Range("A4").Copy($total.Range("A5"))

# Shift rows 2..4 down to 3..5 (captured above before overwriting).
$total.Range("A5").Value = 3
$total.Range("B5").Value = $row4B
$total.Range("C5").Value = $row4C
$total.Range("D5").Value = $row4D

$total.Range("A4").Value = 2
$total.Range("B4").Value = $row3B
$total.Range("C4").Value = $row3C
$total.Range("D4").Value = $row3D

$total.Range("A3").Value = 1
$total.Range("B3").Value = $row2B
$total.Range("C3").Value = $row2C
$total.Range("D3").Value = $row2D

# New top row for 2022-Q1.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# ---------------------------------------------------------------------
# 3. Keep the originally-active sheet selected (adding a sheet makes it
#    active by default).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
